# Update crypto price/volume data cells per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.052.88'

$ws.Range("E2").Value = '  +1.34%  '

$ws.Range("D3").Value = '3.176.86'

$ws.Range("E3").Value = '  +3.66%  '

$ws.Range("E4").Value = '  -0.01%  '

$s = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.81'
$ws.Range("D5").Style = $s

$ws.Range("E5").Value = '  +2.52%  '

$s = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.34'
$ws.Range("D6").Style = $s

$ws.Range("E6").Value = '  +4.64%  '

$s = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = $s

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").Value = '3.169.84'

$ws.Range("E8").Value = '  +3.44%  '

$ws.Range("E9").Value = '  +3.10%  '

$ws.Range("E10").Value = '  +4.47%  '

$ws.Range("E11").Value = '  +1.55%  '

$ws.Range("E12").Value = '  +3.80%  '

$ws.Range("E13").Value = '  +18.81%  '

$ws.Range("E14").Value = '  +7.04%  '

$ws.Range("D15").Value = '3.694.17'

$ws.Range("E15").Value = '  +3.62%  '

$ws.Range("D16").Value = '65.148.20'

$ws.Range("E16").Value = '  +1.42%  '

$ws.Range("B17").Value = 'Polkadot'

$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'

$s = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.20'
$ws.Range("D17").Style = $s

$ws.Range("E17").Value = '  +6.31%  '

$ws.Range("B18").Value = 'WrappedEther'

$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'

$ws.Range("D18").Value = '3.164.91'

$ws.Range("E18").Value = '  +3.24%  '

$ws.Range("E19").Value = '  +1.03%  '

$s = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '513.78'
$ws.Range("D20").Style = $s

$ws.Range("E20").Value = '  +7.18%  '

$s = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.91'
$ws.Range("D21").Style = $s

$ws.Range("E21").Value = '  +6.67%  '

$ws.Range("E22").Value = '  +7.71%  '

$s = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.38'
$ws.Range("D23").Style = $s

$ws.Range("E23").Value = '  +6.99%  '

$ws.Range("E24").Value = '  +3.02%  '

$s = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.35'
$ws.Range("D25").Style = $s

$ws.Range("E25").Value = '  +3.71%  '

$ws.Range("E26").Value = '  +0.04%  '

$s = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.08'
$ws.Range("D27").Style = $s

$ws.Range("E27").Value = '  +11.82%  '

$ws.Range("E28").Value = '  +3.14%  '

$ws.Range("E29").Value = '  +7.65%  '

$s = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '28.13'
$ws.Range("D30").Style = $s

$ws.Range("E30").Value = '  +6.60%  '

$s = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.77'
$ws.Range("D31").Style = $s

$ws.Range("E31").Value = '  +13.13%  '

$ws.Range("E32").Value = '  +6.60%  '

$ws.Range("E33").Value = '  -0.03%  '

$s = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.32'
$ws.Range("D34").Style = $s

$ws.Range("E34").Value = '  +10.31%  '

$ws.Range("E35").Value = '  +6.82%  '

$s = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.70'
$ws.Range("D36").Style = $s

$ws.Range("E36").Value = '  +1.32%  '

$s = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0896'
$ws.Range("D37").Style = $s

$ws.Range("E37").Value = '  +9.98%  '

$s = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '478.71'
$ws.Range("D38").Style = $s

$ws.Range("E38").Value = '  +7.24%  '

$ws.Range("E39").Value = '  +8.64%  '

$ws.Range("E40").Value = '  +3.21%  '

$ws.Range("D41").Value = '3.102.59'

$ws.Range("E41").Value = '  +3.26%  '

$s = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.64'
$ws.Range("D42").Style = $s

$ws.Range("E42").Value = '  +4.55%  '

$ws.Range("E43").Value = '  +3.67%  '

$ws.Range("E44").Value = '  +9.40%  '

$s = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.45'
$ws.Range("D45").Style = $s

$ws.Range("E45").Value = '  +12.40%  '

$s = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '29.42'
$ws.Range("D46").Style = $s

$ws.Range("E46").Value = '  +5.48%  '

$ws.Range("D47").Value = '0.0₃0596'

$ws.Range("E47").Value = '  +14.56%  '

$ws.Range("E48").Value = '  -0.10%  '

$ws.Range("E49").Value = '  +2.15%  '

$ws.Range("E50").Value = '  +9.96%  '

$s = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '121.41'
$ws.Range("D51").Style = $s

$ws.Range("E51").Value = '  +2.33%  '
